$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column H, matching the formatting of the existing header row (B1..G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# "Save" flag values for rows 2..22
$saveValues = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 0
    6  = 0
    7  = 1
    8  = 0
    9  = 1
    10 = 1
    11 = 0
    12 = 0
    13 = 1
    14 = 1
    15 = 0
    16 = 0
    17 = 0
    18 = 0
    19 = 1
    20 = 0
    21 = 1
    22 = 0
}

foreach ($row in $saveValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
